$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 14 (bug #11) was an empty placeholder row. Fill it in with a new bug report,
# re-using the formatting of row 10 (another "严重"-less row with the same
# 是/否BUG = "是" and 状态 = "未解决" values) so the new row inherits the correct
# cell styles automatically.
$ws.Range("D10:M10").Copy()
$ws.Range("D14").PasteSpecial()

# Fill in the actual data for the new bug report.
$ws.Range("D14").Value = 43636
$ws.Range("E14").Value = "沈杰"
$ws.Range("F14").Value = "候天瑞"
$ws.Range("G14").Value = "只连接手环不连接机械手的情况下做动作，动作识别出来但是页面不显示对应图片"
$ws.Range("H14").Value = "未解决"
$ws.Range("I14").Value = "严重"
$ws.Range("J14").Value = "是"
$ws.Range("K14").Value = ""
$ws.Range("L14").Value = ""
$ws.Range("M14").Value = ""

# The severity ("问题等级") column uses a bold font whose color depends on the
# value: red for "严重". Match the formatting used elsewhere in the sheet for
# that value (e.g. row 4's "严重" cell).
$ws.Range("I14").Font.Color = $ws.Range("I4").Font.Color
$ws.Range("I14").Font.Bold = $ws.Range("I4").Font.Bold

# The row grows taller to fit the wrapped description text.
$ws.Rows("14:14").RowHeight = 56

# Update the viewport/selection to reflect where the user ended up after adding
# the new row.
$ws.Range("G14").Select()
